$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.138156
$ws.Range("H2").Value = 0.414468
$ws.Range("I2").Value = 0.0003010053794496939
$ws.Range("J2").Value = 0.0003010053794496939
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.196431
$ws.Range("N2").Value = 0.589293
$ws.Range("O2").Value = 0.09717285149889213
$ws.Range("P2").Value = 0.09717285149889213
$ws.Range("Q2").Value = 0.027138121236
$ws.Range("R2").Value = 0.244243091124
$ws.Range("S2").Value = 2.924955103763278 / 100000
$ws.Range("T2").Value = 2.924955103763278 / 100000

# Row 3
$ws.Range("G3").Value = 0.138156
$ws.Range("H3").Value = 0.414468
$ws.Range("I3").Value = 0.0003010053794496939
$ws.Range("J3").Value = 0.0003010053794496939
$ws.Range("M3").Value = 0.4307096666666667
$ws.Range("O3").Value = 0.2130686423127578
$ws.Range("P3").Value = 0.2130686423127578
$ws.Range("Q3").Value = 0.05950512470800001
$ws.Range("R3").Value = 0.5355461223720001
$ws.Range("S3").Value = 6.413480752818277 / 100000
$ws.Range("T3").Value = 6.413480752818277 / 100000

# Row 4
$ws.Range("G4").Value = 0.138156
$ws.Range("H4").Value = 0.414468
$ws.Range("I4").Value = 0.0003010053794496939
$ws.Range("J4").Value = 0.0003010053794496939
$ws.Range("O4").Value = 0.68975850618835
$ws.Range("P4").Value = 0.68975850618835
$ws.Range("Q4").Value = 0.192633535764
$ws.Range("R4").Value = 1.733701821876
$ws.Range("S4").Value = 0.0002076210208838783
$ws.Range("T4").Value = 0.0002076210208838783

# Row 5
$ws.Range("I5").Value = 0.9878623917146768
$ws.Range("J5").Value = 0.9878623917146769
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.196431
$ws.Range("N5").Value = 0.589293
$ws.Range("O5").Value = 0.09717285149889213
$ws.Range("P5").Value = 0.09717285149889213
$ws.Range("Q5").Value = 89.06395427168198
$ws.Range("R5").Value = 801.5755884451379
$ws.Range("S5").Value = 0.09599340549143069
$ws.Range("T5").Value = 0.09599340549143071

# Row 6
$ws.Range("I6").Value = 0.9878623917146768
$ws.Range("J6").Value = 0.9878623917146769
$ws.Range("M6").Value = 0.4307096666666667
$ws.Range("O6").Value = 0.2130686423127578
$ws.Range("P6").Value = 0.2130686423127578
$ws.Range("S6").Value = 0.2104824985944799
$ws.Range("T6").Value = 0.21048249859448

# Row 7
$ws.Range("I7").Value = 0.9878623917146768
$ws.Range("J7").Value = 0.9878623917146769
$ws.Range("O7").Value = 0.68975850618835
$ws.Range("P7").Value = 0.68975850618835
$ws.Range("S7").Value = 0.6813864876287662
$ws.Range("T7").Value = 0.6813864876287662

# Row 8
$ws.Range("G8").Value = 5.432785666666668
$ws.Range("I8").Value = 0.01183660290587349
$ws.Range("J8").Value = 0.01183660290587349
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.196431
$ws.Range("N8").Value = 0.589293
$ws.Range("O8").Value = 0.09717285149889213
$ws.Range("P8").Value = 0.09717285149889213
$ws.Range("Q8").Value = 1.067167521289
$ws.Range("R8").Value = 9.604507691601
$ws.Range("S8").Value = 0.0011501964564238
$ws.Range("T8").Value = 0.0011501964564238

# Row 9
$ws.Range("G9").Value = 5.432785666666668
$ws.Range("I9").Value = 0.01183660290587349
$ws.Range("J9").Value = 0.01183660290587349
$ws.Range("M9").Value = 0.4307096666666667
$ws.Range("Q9").Value = 2.339953303561445
$ws.Range("R9").Value = 21.05957973205301
$ws.Range("S9").Value = 0.002522008910749709
$ws.Range("T9").Value = 0.002522008910749709

# Row 10
$ws.Range("G10").Value = 5.432785666666668
$ws.Range("I10").Value = 0.01183660290587349
$ws.Range("J10").Value = 0.01183660290587349
$ws.Range("O10").Value = 0.68975850618835
$ws.Range("P10").Value = 0.68975850618835
$ws.Range("Q10").Value = 7.575036277961002
$ws.Range("R10").Value = 68.17532650164901
$ws.Range("S10").Value = 0.008164397538699984
$ws.Range("T10").Value = 0.008164397538699984
